# Backend listo para despliegue
# Adds a new data row (row 2) with a registration record, widens column D
# to match the other "20"-wide columns, and refreshes the page setup /
# print settings for the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New record (row 2) -----------------------------------------------
$ws.Range("A2").Value = "15/10/2025, 12:17:46 a. m."
$ws.Range("B2").Value = "Dania"
$ws.Range("C2").Value = "Soto"
$ws.Range("D2").Value = "Voleibol"
$ws.Range("E2").Value = "female"
$ws.Range("F2").Value = "Jutiapa"
$ws.Range("G2").Value = "No"
$ws.Range("H2").Value = "Toyota"

# --- Column widths -------------------------------------------------------
# Column D ("Deporte") now matches the 20-wide columns A:C instead of
# sharing the 15-wide formatting with column E.
$ws.Columns.Item(4).ColumnWidth = 19.17

# --- Page setup / print options -----------------------------------------
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.FirstPageNumber = 1
$ws.PageSetup.UseFirstPageNumber = $true
$ws.PageSetup.Copies = 1
